# -----------------------------------------------------------------------
# Updates the 'cryptos' price/volume snapshot table (columns D = Price,
# E = Volume(1h)); a few rows also get new Coin/Link values where two
# coins swapped rank position between snapshots (rows 43/44, 47/48).
#
# Price values such as '232.51' or '1.00' look like numbers to Excel, so
# a leading apostrophe (quote-prefix) is used to force them to stay text,
# matching the column's existing text formatting - exactly like a user
# typing '232.51 into the cell.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '37.741.38'
$ws.Range('E2').Value = '  +1.30%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.091.13'

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.07%  '

# Row 5: BNB
$ws.Range('D5').Value = '''232.51'
$ws.Range('E5').Value = '  -0.11%  '

# Row 6: XRP
$ws.Range('E6').Value = '  +0.49%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.05%  '

# Row 8: Solana
$ws.Range('D8').Value = '''57.77'
$ws.Range('E8').Value = '  +2.35%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  +1.91%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +2.36%  '

# Row 11: TRON
$ws.Range('E11').Value = '  +3.05%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '2.387.80'
$ws.Range('E12').Value = '  +1.05%  '

# Row 13: Chainlink
$ws.Range('D13').Value = '''14.43'
$ws.Range('E13').Value = '  -0.98%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '''21.11'
$ws.Range('E14').Value = '  +2.55%  '

# Row 15: Polygon
$ws.Range('D15').Value = '''0.766'
$ws.Range('E15').Value = '  -1.07%  '

# Row 16: Polkadot
$ws.Range('E16').Value = '  +2.40%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.085.13'
$ws.Range('E17').Value = '  +1.16%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '37.668.51'
$ws.Range('E18').Value = '  +1.23%  '

# Row 19: Uniswap
$ws.Range('D19').Value = '''6.13'
$ws.Range('E19').Value = '  -2.65%  '

# Row 20: Litecoin
$ws.Range('D20').Value = '''70.60'
$ws.Range('E20').Value = '  +2.03%  '

# Row 21: ShibaInu
$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('E21').Value = '  +1.63%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''228.03'
$ws.Range('E22').Value = '  +0.99%  '

# Row 23: Dai
$ws.Range('E23').Value = '  -0.02%  '

# Row 24: Toncoin
$ws.Range('E24').Value = '  -1.16%  '

# Row 25: PancakeSwap
$ws.Range('D25').Value = '''2.37'
$ws.Range('E25').Value = '  -0.18%  '

# Row 26: Monero
$ws.Range('D26').Value = '''167.89'
$ws.Range('E26').Value = '  +0.56%  '

# Row 27: Kaspa
$ws.Range('E27').Value = '  +9.26%  '

# Row 28: Cosmos
$ws.Range('D28').Value = '''8.95'
$ws.Range('E28').Value = '  +2.35%  '

# Row 29: ImmutableX
$ws.Range('E29').Value = '  -0.84%  '

# Row 30: EthereumClassic
$ws.Range('E30').Value = '  +2.44%  '

# Row 31: Stellar
$ws.Range('E31').Value = '  +1.33%  '

# Row 32: Filecoin
$ws.Range('E32').Value = '  +4.27%  '

# Row 33: Hedera
$ws.Range('E33').Value = '  +1.74%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').Value = '''4.57'
$ws.Range('E34').Value = '  -0.07%  '

# Row 35: LidoDAOToken
$ws.Range('D35').Value = '''2.50'
$ws.Range('E35').Value = '  +0.19%  '

# Row 36: WEMIXToken
$ws.Range('E36').Value = '  +4.31%  '

# Row 37: RenderToken
$ws.Range('D37').Value = '''3.40'
$ws.Range('E37').Value = '  +5.34%  '

# Row 38: BinanceUSD
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  -0.10%  '

# Row 39: THORChain
$ws.Range('D39').Value = '''5.39'
$ws.Range('E39').Value = '  -5.01%  '

# Row 40: Cronos
$ws.Range('D40').Value = '''0.0993'
$ws.Range('E40').Value = '  +5.93%  '

# Row 41: HuobiToken
$ws.Range('E41').Value = '  -0.46%  '

# Row 42: Aave
$ws.Range('D42').Value = '''97.87'
$ws.Range('E42').Value = '  +2.25%  '

# Row 43: Maker
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0214'
$ws.Range('E43').Value = '  +0.76%  '

# Row 44: VeChain
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.453.87'
$ws.Range('E44').Value = '  -1.01%  '

# Row 45: TrustWalletToken
$ws.Range('E45').Value = '  -0.33%  '

# Row 46: ARBITRUM
$ws.Range('E46').Value = '  +3.89%  '

# Row 47: InjectiveProtocol
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = '''4.06'
$ws.Range('E47').Value = '  -5.35%  '

# Row 48: FTXToken
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''15.63'
$ws.Range('E48').Value = '  +3.93%  '

# Row 49: FraxShare
$ws.Range('E49').Value = '  +3.08%  '

# Row 50: MXToken
$ws.Range('E50').Value = '  +2.04%  '

# Row 51: RocketPoolETH
$ws.Range('D51').Value = '2.283.88'
$ws.Range('E51').Value = '  +1.47%  '
